# Apply "Add data for 2022-10-17" update:
#  - Rename sheet to reflect new "through" date (10-08 -> 10-09)
#  - Update October row label text (10-08 -> 10-09)
#  - Update October row (row 11) and Total row (row 12) figures for
#    columns B..I (years 2015..2022)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "Through 2022-10-09"

# Update the October row label (shared string used by A11)
$ws.Range("A11").Value = "October (through 10-09)"

# Update October row (row 11) values
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 16
$ws.Range("D11").Value = 17
$ws.Range("E11").Value = 23
$ws.Range("F11").Value = 8
$ws.Range("G11").Value = 38
$ws.Range("H11").Value = 61
$ws.Range("I11").Value = 30

# Update Total row (row 12) values
$ws.Range("B12").Value = 233
$ws.Range("C12").Value = 445
$ws.Range("D12").Value = 644
$ws.Range("E12").Value = 571
$ws.Range("F12").Value = 430
$ws.Range("G12").Value = 939
$ws.Range("H12").Value = 1308
$ws.Range("I12").Value = 1308

$wb.Save()
